# Apply updated crypto price/volume data (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.524.87'
$ws.Range('E2').Value = '  +0.45%  '

$ws.Range('D3').Value = '1.639.32'
$ws.Range('E3').Value = '  -0.72%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.72'
$ws.Range('E5').Value = '  -0.43%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.532'
$ws.Range('E6').Value = '  +4.62%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.01'
$ws.Range('E8').Value = '  -5.06%  '

$ws.Range('E9').Value = '  -2.12%  '

$ws.Range('E10').Value = '  -0.73%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0890'

$ws.Range('D12').Value = '1.869.14'
$ws.Range('E12').Value = '  -0.85%  '

$ws.Range('D13').Value = '1.639.96'
$ws.Range('E13').Value = '  -0.63%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.03'
$ws.Range('E14').Value = '  -1.41%  '

$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.564'
$ws.Range('E15').Value = '  -1.83%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.35'
$ws.Range('E16').Value = '  -2.49%  '

$ws.Range('D17').Value = '27.486.13'
$ws.Range('E17').Value = '  +0.29%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.92'
$ws.Range('E18').Value = '  -1.91%  '

$ws.Range('E19').Value = '  +2.99%  '

$ws.Range('D20').Value = '0.0₃0724'
$ws.Range('E20').Value = '  -0.33%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  -0.01%  '

$ws.Range('E22').Value = '  -1.80%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.97'
$ws.Range('E23').Value = '  +7.09%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.94'
$ws.Range('E24').Value = '  -3.56%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.52'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.97'
$ws.Range('E26').Value = '  -3.29%  '

$ws.Range('E27').Value = '  +1.61%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.07%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.60'
$ws.Range('E29').Value = '  -3.02%  '

$ws.Range('E30').Value = '  -0.28%  '

$ws.Range('E31').Value = '  -1.93%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.29'
$ws.Range('E32').Value = '  -0.51%  '

$ws.Range('E33').Value = '  +1.97%  '

$ws.Range('D34').Value = '1.422.50'
$ws.Range('E34').Value = '  -2.55%  '

$ws.Range('E35').Value = '  +2.32%  '

$ws.Range('E36').Value = '  -1.86%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.571'
$ws.Range('E37').Value = '  -0.44%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.879'
$ws.Range('E38').Value = '  -3.46%  '

$ws.Range('E39').Value = '  -1.81%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.891'
$ws.Range('E40').Value = '  +13.52%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.04'
$ws.Range('E41').Value = '  -0.15%  '

$ws.Range('E42').Value = '  +0.04%  '

$ws.Range('E43').Value = '  -0.28%  '

$ws.Range('E44').Value = '  +1.59%  '

$ws.Range('E45').Value = '  +1.48%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.98'
$ws.Range('E46').Value = '  -0.80%  '

$ws.Range('D47').Value = '1.779.40'
$ws.Range('E47').Value = '  -0.80%  '

$ws.Range('E48').Value = '  -3.09%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '86.13'
$ws.Range('E49').Value = '  -2.75%  '

$ws.Range('D50').Value = '0.0₆0103'
$ws.Range('E50').Value = '  -2.39%  '

$ws.Range('E51').Value = '  -2.04%  '
